$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '42.614.69'
$ws.Cells.Item(2, 5).Value = '  -0.02%  '

$ws.Cells.Item(3, 4).Value = '2.292.85'
$ws.Cells.Item(3, 5).Value = '  +0.60%  '

$ws.Cells.Item(4, 5).Value = '  +0.05%  '

$ws.Cells.Item(5, 4).Value = '''301.27'
$ws.Cells.Item(5, 5).Value = '  -1.16%  '

$ws.Cells.Item(6, 4).Value = '''95.49'
$ws.Cells.Item(6, 5).Value = '  -0.98%  '

$ws.Cells.Item(7, 5).Value = '  -0.10%  '

$ws.Cells.Item(8, 5).Value = '  +0.12%  '

$ws.Cells.Item(9, 5).Value = '  -1.90%  '

$ws.Cells.Item(10, 4).Value = '''34.31'
$ws.Cells.Item(10, 5).Value = '  -3.39%  '

$ws.Cells.Item(11, 4).Value = '''18.96'
$ws.Cells.Item(11, 5).Value = '  +4.22%  '

$ws.Cells.Item(12, 5).Value = '  -0.64%  '

$ws.Cells.Item(13, 5).Value = '  +0.09%  '

$ws.Cells.Item(14, 5).Value = '  +0.00%  '

$ws.Cells.Item(15, 4).Value = '2.653.23'
$ws.Cells.Item(15, 5).Value = '  +0.74%  '

$ws.Cells.Item(16, 4).Value = '2.291.33'
$ws.Cells.Item(16, 5).Value = '  +0.38%  '

$ws.Cells.Item(17, 4).Value = '''0.778'
$ws.Cells.Item(17, 5).Value = '  +0.22%  '

$ws.Cells.Item(18, 4).Value = '42.549.56'
$ws.Cells.Item(18, 5).Value = '  -0.06%  '

$ws.Cells.Item(19, 4).Value = '''12.14'
$ws.Cells.Item(19, 5).Value = '  -6.52%  '

$ws.Cells.Item(20, 5).Value = '  -0.69%  '

$ws.Cells.Item(21, 5).Value = '  -0.43%  '

$ws.Cells.Item(22, 4).Value = '''67.68'
$ws.Cells.Item(22, 5).Value = '  +0.80%  '

$ws.Cells.Item(24, 4).Value = '''235.12'
$ws.Cells.Item(24, 5).Value = '  -0.21%  '

$ws.Cells.Item(25, 5).Value = '  +0.08%  '

$ws.Cells.Item(26, 4).Value = '''2.41'
$ws.Cells.Item(26, 5).Value = '  -1.66%  '

$ws.Cells.Item(27, 5).Value = '  -3.41%  '

$ws.Cells.Item(28, 5).Value = '  +14.40%  '

$ws.Cells.Item(29, 4).Value = '''165.68'

$ws.Cells.Item(30, 5).Value = '  -0.26%  '

$ws.Cells.Item(31, 4).Value = '''31.75'
$ws.Cells.Item(31, 5).Value = '  -3.88%  '

$ws.Cells.Item(32, 5).Value = '  +0.02%  '

$ws.Cells.Item(33, 5).Value = '  +0.21%  '

$ws.Cells.Item(34, 4).Value = '''17.46'
$ws.Cells.Item(34, 5).Value = '  -0.81%  '

$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 4).Value = '''0.0697'
$ws.Cells.Item(35, 5).Value = '  +1.05%  '

$ws.Cells.Item(36, 2).Value = 'RenderToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(36, 4).Value = '''4.40'
$ws.Cells.Item(36, 5).Value = '  -7.19%  '

$ws.Cells.Item(37, 5).Value = '  -2.60%  '

$ws.Cells.Item(38, 5).Value = '  -1.17%  '

$ws.Cells.Item(39, 5).Value = '  -0.02%  '

$ws.Cells.Item(40, 5).Value = '  -1.37%  '

$ws.Cells.Item(41, 4).Value = '''2.67'
$ws.Cells.Item(41, 5).Value = '  -0.63%  '

$ws.Cells.Item(42, 4).Value = '''20.14'
$ws.Cells.Item(42, 5).Value = '  +11.76%  '

$ws.Cells.Item(43, 4).Value = '1.962.28'
$ws.Cells.Item(43, 5).Value = '  -2.09%  '

$ws.Cells.Item(44, 4).Value = '''10.42'
$ws.Cells.Item(44, 5).Value = '  +4.44%  '

$ws.Cells.Item(45, 5).Value = '  -0.14%  '

$ws.Cells.Item(46, 5).Value = '  -0.43%  '

$ws.Cells.Item(47, 5).Value = '  -0.82%  '

$ws.Cells.Item(48, 4).Value = '2.519.00'
$ws.Cells.Item(48, 5).Value = '  +0.68%  '

$ws.Cells.Item(49, 2).Value = 'MultiversX'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(49, 4).Value = '''52.98'
$ws.Cells.Item(49, 5).Value = '  -1.03%  '

$ws.Cells.Item(50, 2).Value = 'HuobiToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(50, 4).Value = '''2.79'
$ws.Cells.Item(50, 5).Value = '  -1.78%  '

$ws.Cells.Item(51, 4).Value = '''71.16'
$ws.Cells.Item(51, 5).Value = '  -0.13%  '
